$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells keep their text representation (values like
# "412.06" or "3.415.24" would otherwise be auto-converted to numbers/dates by Excel).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.124.49"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.415.24"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.06"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.32"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -2.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.68"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.954.66"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.44"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.430.01"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.75"
$ws.Range("E18").Value = "  +5.19%  "
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.149.99"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "472.17"
$ws.Range("E21").Value = "  +6.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.58"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("E23").Value = "  +3.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.04"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.30"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.79"
$ws.Range("E26").Value = "  +11.23%  "
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.75"
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.91"
$ws.Range("E34").Value = "  -4.67%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.06"
$ws.Range("E36").Value = "  +8.29%  "
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +3.69%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.70"
$ws.Range("E42").Value = "  +12.64%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.42"
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("E45").Value = "  +4.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.33"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.35"
$ws.Range("E47").Value = "  +18.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.33"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0535"
$ws.Range("E49").Value = "  +25.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.32"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.94"
$ws.Range("E51").Value = "  +7.07%  "
